# Rename account bills' "Date" column header to "Transfer date".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# The "Date" header lives in C1 on the Transactions sheet; rename it.
$ws.Range("C1").Value = "Transfer date"

# Reflect the author's resulting selection (clicked back on the renamed header).
$ws.Activate()
$ws.Range("C1").Select() | Out-Null
